$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @(0.75, 4236727.498161, 21243.694582, 4215483.803579, 13991.442363, 1760426.373385, 14071.290441, 1746355.082944, 15981.531147),
        @(0.8, 5117794.229655, 21393.73706, 5096400.492595, 14391.29253366667, 2050415.88006, 13914.390853, 2036501.489207, 16453.033397),
        @(0.85, 5880494.184356, 21702.176071, 5858792.008285, 14916.40278166667, 2266081.715755, 13768.710049, 2252313.005706, 17045.34805033333),
        @(0.9, 6426087.193499, 21818.546042, 6404268.647457, 15350.98245466666, 2395551.12106, 13640.013058, 2381911.108002, 17559.29677233333),
        @(0.95, 6841596.896981, 22087.582224, 6819509.314757, 15900.36611133333, 2470749.551649, 13579.795378, 2457169.756271, 18182.396295)
)

$startRow = 18
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $vals[$j]
    }
}
